$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last tab (3rd sheet)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ansp_lists"

# Header row
$ws.Range("A1").Value = "ANSP_CODE"
$ws.Range("B1").Value = "ANSP_NAME"
$ws.Range("C1").Value = "PRU_ID"

# Data rows (ANSP_CODE, ANSP_NAME, PRU_ID)
$arr = New-Object "object[,]" 41,3
$arr[0,0] = "IE_ANSP"; $arr[0,1] = "AirNav Ireland"; $arr[0,2] = 17
$arr[1,0] = "AL_ANSP"; $arr[1,1] = "Albcontrol"; $arr[1,2] = 24
$arr[2,0] = "CZ_ANSP"; $arr[2,1] = "ANS CR"; $arr[2,2] = 2
$arr[3,0] = "AM_ANSP"; $arr[3,1] = "ARMATS"; $arr[3,2] = 44
$arr[4,0] = "AT_ANSP"; $arr[4,1] = "Austro Control"; $arr[4,2] = 4
$arr[5,0] = "NO_ANSP"; $arr[5,1] = "Avinor"; $arr[5,2] = 50
$arr[6,0] = "BA_ANSP"; $arr[6,1] = "BHANSA"; $arr[6,2] = 45
$arr[7,0] = "BG_ANSP"; $arr[7,1] = "BULATSA"; $arr[7,2] = 3
$arr[8,0] = "HR_ANSP"; $arr[8,1] = "Croatia Control"; $arr[8,2] = 6
$arr[9,0] = "CY_ANSP"; $arr[9,1] = "DCAC Cyprus"; $arr[9,2] = 7
$arr[10,0] = "DE_ANSP"; $arr[10,1] = "DFS"; $arr[10,2] = 8
$arr[11,0] = "TR_ANSP"; $arr[11,1] = "DHMI"; $arr[11,2] = 9
$arr[12,0] = "FR_ANSP"; $arr[12,1] = "DSNA"; $arr[12,2] = 10
$arr[13,0] = "EE_ANSP"; $arr[13,1] = "EANS"; $arr[13,2] = 11
$arr[14,0] = "ES_ANSP"; $arr[14,1] = "ENAIRE"; $arr[14,2] = 1
$arr[15,0] = "IT_ANSP"; $arr[15,1] = "ENAV"; $arr[15,2] = 12
$arr[16,0] = "FI_ANSP"; $arr[16,1] = "Fintraffic ANS"; $arr[16,2] = 13
$arr[17,0] = "GR_ANSP"; $arr[17,1] = "HASP"; $arr[17,2] = 15
$arr[18,0] = "HU_ANSP"; $arr[18,1] = "HungaroControl (EC)"; $arr[18,2] = 16
$arr[19,0] = "IS_ANSP"; $arr[19,1] = "Isavia"; $arr[19,2] = 46
$arr[20,0] = "IL_ANSP"; $arr[20,1] = "Israel AA"; $arr[20,2] = 57
$arr[21,0] = "SE_ANSP"; $arr[21,1] = "LFV"; $arr[21,2] = 33
$arr[22,0] = "LV_ANSP"; $arr[22,1] = "LGS"; $arr[22,2] = 18
$arr[23,0] = "SK_ANSP"; $arr[23,1] = "LPS"; $arr[23,2] = 19
$arr[24,0] = "NL_ANSP"; $arr[24,1] = "LVNL"; $arr[24,2] = 20
$arr[25,0] = "MT_ANSP"; $arr[25,1] = "MATS"; $arr[25,2] = 21
$arr[26,0] = "MK_ANSP"; $arr[26,1] = "M-NAV"; $arr[26,2] = 14
$arr[27,0] = "MD_ANSP"; $arr[27,1] = "MOLDATSA"; $arr[27,2] = 22
$arr[28,0] = "MAS_ANSP"; $arr[28,1] = "MUAC"; $arr[28,2] = 23
$arr[29,0] = "GB_ANSP"; $arr[29,1] = "NATS (Continental)"; $arr[29,2] = 26
$arr[30,0] = "PT_ANSP"; $arr[30,1] = "NAV Portugal (Continental)"; $arr[30,2] = 27
$arr[31,0] = "DK_ANSP"; $arr[31,1] = "NAVIAIR"; $arr[31,2] = 28
$arr[32,0] = "MA_ANSP"; $arr[32,1] = "ONDA"; $arr[32,2] = 56
$arr[33,0] = "LT_ANSP"; $arr[33,1] = "Oro Navigacija"; $arr[33,2] = 29
$arr[34,0] = "PL_ANSP"; $arr[34,1] = "PANSA"; $arr[34,2] = 39
$arr[35,0] = "RO_ANSP"; $arr[35,1] = "ROMATSA"; $arr[35,2] = 30
$arr[36,0] = "GE_ANSP"; $arr[36,1] = "Sakaeronavigatsia"; $arr[36,2] = 53
$arr[37,0] = "BE_ANSP"; $arr[37,1] = "skeyes"; $arr[37,2] = 5
$arr[38,0] = "CH_ANSP"; $arr[38,1] = "Skyguide"; $arr[38,2] = 31
$arr[39,0] = "SI_ANSP"; $arr[39,1] = "Slovenia Control"; $arr[39,2] = 32
$arr[40,0] = "MERS_ANSP"; $arr[40,1] = "SMATSA"; $arr[40,2] = 42

$ws.Range("A2:C42").Value = $arr

# Turn the range into an Excel table ("Table_ANSP_NAMES")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:C42"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table_ANSP_NAMES"

# Restore the selection / active cell on the new sheet
$ws.Range("B7").Select()

Write-Output "ansp_lists sheet created"
